$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, pushing the existing rows 40-103 down to 41-104.
$ws.Rows("40").Insert()

# Populate the newly inserted row 40 with the new record's data.
$ws.Range("A40").Value = 11
$ws.Range("B40").Value = "Vega Monumental Concepción"
$ws.Range("C40").Value = "Bíobío"
$ws.Range("D40").Value = 44771
$ws.Range("E40").Value = 8
$ws.Range("F40").Value = 100112001
$ws.Range("G40").Value = "Berenjena"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 100
$ws.Range("K40").Value = 11000
$ws.Range("L40").Value = 12000
$ws.Range("M40").Value = 11500
$ws.Range("N40").Value = "$/caja 60 unidades"
$ws.Range("O40").Value = "Región de Arica y Parinacota"
$ws.Range("P40").Value = 192
$ws.Range("Q40").Value = 60
$ws.Range("R40").Value = "Hortaliza"
